# Pavani Gandepalli daily status on Feb 1 2021
# Adds a new dated status row (row 23) to the B&I status sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 23: new daily status entry ---

# Copy formatting from existing rows into row 23 so the new row matches the
# look of the rest of the table (date style for col A, wrap-text style for
# the text columns B and D).
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A23").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B9").Copy() | Out-Null
$ws.Range("B23").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B9").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Application.CutCopyMode = $false

$ws.Cells.Item(23, 1).Value = 44229
$ws.Cells.Item(23, 4).Value = "FindDigits.txt`r`nFunctionsInC.txt"
$ws.Cells.Item(23, 2).Value = "1. I have completed bitwise operators 1 program, pointers 5 programs from given list of programs and pushed to github`r`n2. I worked on performance testing with linpack app, Vellamo app `r`n3. I have worked on camera few Camera test cases and tested them on my mobile and recored the results`r`n4. Attended Srinivas session about test cases writing and validation of my `r`n5. I have completed 2 hacker rank programs today "

$ws.Rows.Item(23).RowHeight = 180

$ws.Range("C23").Select() | Out-Null
